$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - 常数项 (Constant)
$ws.Range("B2").Value = "9.4695***"
$ws.Range("C2").Value = "'0.0278"
$ws.Range("D2").Value = "14.3785***"
$ws.Range("E2").Value = "'1.9272"
$ws.Range("F2").Value = "'7.4610"

# Row 3 - did
$ws.Range("B3").Value = "'0.0226"
$ws.Range("C3").Value = "'0.0288"
$ws.Range("D3").Value = "'0.0275"
$ws.Range("E3").Value = "'0.0282"
$ws.Range("F3").Value = "'0.9761"
$ws.Range("G3").Value = "'0.3291"

# Row 4 - ln_pgdp
$ws.Range("D4").Value = "-0.3598***"
$ws.Range("E4").Value = "'0.1333"
$ws.Range("F4").Value = "'-2.6986"
$ws.Range("G4").Value = "'0.0070"

# Row 5 - ln_pop_density
$ws.Range("D5").Value = "'-0.2505"
$ws.Range("E5").Value = "'0.1670"
$ws.Range("F5").Value = "'-1.4999"
$ws.Range("G5").Value = "'0.1337"

# Row 6 - industrial_advanced
$ws.Range("D6").Value = "-0.0841***"
$ws.Range("E6").Value = "'0.0243"
$ws.Range("F6").Value = "'-3.4612"
$ws.Range("G6").Value = "'0.0005"

# Row 7 - fdi_openness
$ws.Range("D7").Value = "'-0.9608"
$ws.Range("E7").Value = "'0.7121"
$ws.Range("F7").Value = "'-1.3493"
$ws.Range("G7").Value = "'0.1773"

# Row 8 - ln_road_area
$ws.Range("D8").Value = "0.1150***"
$ws.Range("E8").Value = "'0.0445"
$ws.Range("F8").Value = "'2.5823"
$ws.Range("G8").Value = "'0.0099"

# Row 9 - 样本量 (Sample size) - numeric values
$ws.Range("B9").Value = 2846
$ws.Range("D9").Value = 2846

# Row 11 - R2
$ws.Range("B11").Value = "'0.9677"
$ws.Range("D11").Value = "'0.9697"

# Row 12 - 调整R2 (Adjusted R2)
$ws.Range("B12").Value = "'0.9651"
$ws.Range("D12").Value = "'0.9671"
